# Update the "Online Test - Export Question" sheet's header row:
# the "Repeat Time" column is renamed to "Times of listening" (column G / header row 5),
# the column is widened to fit the new label, and the view is scrolled down
# to the data rows with a new active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SampleTestSet")

# Rename the header text in G5 from "Repeat Time" to "Times of listening".
$ws.Range("G5").Value = "Times of listening"

# Widen column G so the longer label fits (was ~13.14 chars, now ~17.43 chars).
$ws.Columns.Item(7).ColumnWidth = 16.6

# Update the sheet's scroll position / active selection (was topLeftCell B1 / H10).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 2
$ws.Range("H18").Select()
